$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing B3 note: text expanded with an additional sentence,
# and the row grows taller to fit the extra wrapped text.
$ws.Range("B3").Value = "Jätkake koodi kirjutamist. Lisatud pildid ja joonised tegelastele.                        Tausta pildite teeme esimesed fotod."
$ws.Rows.Item(3).RowHeight = 75

# Append a new work-log entry on row 4.
$ws.Range("A4").Value = "20.05.2024"
$ws.Range("B4").Value = "Uus plaan ellu viia.                     Mängu süžeeplaan on kirjutatud.                                      Dialoogi koostamine.                Alustage tegelaste joonistamist."
$ws.Range("B4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 90

# Move/collapse the selection to F4, matching the saved view state.
[void]$ws.Range("F4").Select()
